$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 0
$ws.Cells.Item(2, 2).Value = "falling"
$ws.Cells.Item(2, 3).Value = -2.243617534637451
$ws.Cells.Item(2, 4).Value = 1.99995231628418
$ws.Cells.Item(2, 5).Value = 0.6877330541610718
$ws.Cells.Item(2, 6).Value = -0.06436660403446501
$ws.Cells.Item(2, 7).Value = 0.02161267820907678
$ws.Cells.Item(2, 8).Value = -0.02475332161006712

$ws.Cells.Item(3, 1).Value = 100
$ws.Cells.Item(3, 2).Value = "falling"
$ws.Cells.Item(3, 3).Value = -1.472963809967041
$ws.Cells.Item(3, 4).Value = 1.918478488922119
$ws.Cells.Item(3, 5).Value = 0.357224702835083
$ws.Cells.Item(3, 6).Value = -0.002483299958681726
$ws.Cells.Item(3, 7).Value = -0.07227465348399185
$ws.Cells.Item(3, 8).Value = -0.05542938006312945

$ws.Cells.Item(4, 1).Value = 200
$ws.Cells.Item(4, 2).Value = "falling"
$ws.Cells.Item(4, 3).Value = -2.048044681549072
$ws.Cells.Item(4, 4).Value = 1.973574161529541
$ws.Cells.Item(4, 5).Value = 0.7422993183135986
$ws.Cells.Item(4, 6).Value = 0.05864306564604772
$ws.Cells.Item(4, 7).Value = -0.01641368349690153
$ws.Cells.Item(4, 8).Value = -0.03955352419744367

$ws.Cells.Item(5, 1).Value = 300
$ws.Cells.Item(5, 2).Value = "falling"
$ws.Cells.Item(5, 3).Value = -1.858012676239014
$ws.Cells.Item(5, 4).Value = 2.016510963439941
$ws.Cells.Item(5, 5).Value = 0.136701762676239
$ws.Cells.Item(5, 6).Value = 0.0103050321340562
$ws.Cells.Item(5, 7).Value = 0.04435412786172135
$ws.Cells.Item(5, 8).Value = 0.04177123000440386

$ws.Cells.Item(6, 1).Value = 400
$ws.Cells.Item(6, 2).Value = "falling"
$ws.Cells.Item(6, 3).Value = -1.831352233886719
$ws.Cells.Item(6, 4).Value = 1.850147724151612
$ws.Cells.Item(6, 5).Value = 0.5375880002975464
$ws.Cells.Item(6, 6).Value = 0.02360463045213534
$ws.Cells.Item(6, 7).Value = 0.01676559415848359
$ws.Cells.Item(6, 8).Value = 0.06694286187057909

$ws.Cells.Item(7, 1).Value = 500
$ws.Cells.Item(7, 2).Value = "falling"
$ws.Cells.Item(7, 3).Value = -1.778313636779785
$ws.Cells.Item(7, 4).Value = 1.832831382751465
$ws.Cells.Item(7, 5).Value = 0.6446512937545776
$ws.Cells.Item(7, 6).Value = 0.02699094848788307
$ws.Cells.Item(7, 7).Value = -0.01045774781833528
$ws.Cells.Item(7, 8).Value = 0.02604145086977801

$ws.Cells.Item(8, 1).Value = 600
$ws.Cells.Item(8, 2).Value = "falling"
$ws.Cells.Item(8, 3).Value = -1.913262844085693
$ws.Cells.Item(8, 4).Value = 1.87591552734375
$ws.Cells.Item(8, 5).Value = 0.3856300711631775
$ws.Cells.Item(8, 6).Value = -0.01472052470173514
$ws.Cells.Item(8, 7).Value = 0.001075654007170485
$ws.Cells.Item(8, 8).Value = -0.02027806955511146

$ws.Cells.Item(9, 1).Value = 700
$ws.Cells.Item(9, 2).Value = "falling"
$ws.Cells.Item(9, 3).Value = -2.06378698348999
$ws.Cells.Item(9, 4).Value = 1.738418102264404
$ws.Cells.Item(9, 5).Value = 0.8020429015159607
$ws.Cells.Item(9, 6).Value = -0.0279470856260994
$ws.Cells.Item(9, 7).Value = -0.001367806902398176
$ws.Cells.Item(9, 8).Value = -0.02752877472211478

$ws.Cells.Item(10, 1).Value = 800
$ws.Cells.Item(10, 2).Value = "falling"
$ws.Cells.Item(10, 3).Value = -1.788212299346924
$ws.Cells.Item(10, 4).Value = 1.889766216278076
$ws.Cells.Item(10, 5).Value = 0.3115454912185669
$ws.Cells.Item(10, 6).Value = -0.02199114867202608
$ws.Cells.Item(10, 7).Value = 0.004966600231178349
$ws.Cells.Item(10, 8).Value = -0.03471972080676452

$ws.Cells.Item(11, 1).Value = 900
$ws.Cells.Item(11, 2).Value = "falling"
$ws.Cells.Item(11, 3).Value = -1.937233924865723
$ws.Cells.Item(11, 4).Value = 1.803182125091553
$ws.Cells.Item(11, 5).Value = 0.8862782716751099
$ws.Cells.Item(11, 6).Value = -0.01893682202891156
$ws.Cells.Item(11, 7).Value = 0.0004515091119253129
$ws.Cells.Item(11, 8).Value = 0.08639095043358573

$ws.Cells.Item(12, 1).Value = 1000
$ws.Cells.Item(12, 2).Value = "falling"
$ws.Cells.Item(12, 3).Value = -2.046533107757568
$ws.Cells.Item(12, 4).Value = 1.845665454864502
$ws.Cells.Item(12, 5).Value = 0.7709857225418091
$ws.Cells.Item(12, 6).Value = -0.007516298443078927
$ws.Cells.Item(12, 7).Value = 0.01816660165786737
$ws.Cells.Item(12, 8).Value = 0.06252737019373024

$ws.Cells.Item(13, 1).Value = 1100
$ws.Cells.Item(13, 2).Value = "falling"
$ws.Cells.Item(13, 3).Value = -1.904958724975586
$ws.Cells.Item(13, 4).Value = 1.821238040924072
$ws.Cells.Item(13, 5).Value = 0.7353460788726807
$ws.Cells.Item(13, 6).Value = -0.01159979990157086
$ws.Cells.Item(13, 7).Value = 0.006480483672298237
$ws.Cells.Item(13, 8).Value = 0.04596096941310422

$ws.Cells.Item(14, 1).Value = 1200
$ws.Cells.Item(14, 2).Value = "falling"
$ws.Cells.Item(14, 3).Value = -1.893377780914306
$ws.Cells.Item(14, 4).Value = 1.849058628082276
$ws.Cells.Item(14, 5).Value = 0.4389755129814148
$ws.Cells.Item(14, 6).Value = -0.01335271759687552
$ws.Cells.Item(14, 7).Value = 0.001062374369686679
$ws.Cells.Item(14, 8).Value = 0.02341207367894439

$ws.Cells.Item(15, 1).Value = 1300
$ws.Cells.Item(15, 2).Value = "falling"
$ws.Cells.Item(15, 3).Value = -1.822815418243408
$ws.Cells.Item(15, 4).Value = 1.71094799041748
$ws.Cells.Item(15, 5).Value = 0.5070880651473999
$ws.Cells.Item(15, 6).Value = -0.02131388465995362
$ws.Cells.Item(15, 7).Value = -0.005225553865665965
$ws.Cells.Item(15, 8).Value = 0.007270624132259934

$ws.Cells.Item(16, 1).Value = 1400
$ws.Cells.Item(16, 2).Value = "falling"
$ws.Cells.Item(16, 3).Value = -1.894952774047852
$ws.Cells.Item(16, 4).Value = 1.727813243865967
$ws.Cells.Item(16, 5).Value = 0.551978588104248
$ws.Cells.Item(16, 6).Value = -0.03625352340547929
$ws.Cells.Item(16, 7).Value = -0.01064366328975433
$ws.Cells.Item(16, 8).Value = 0.008558753068032388

$ws.Cells.Item(17, 1).Value = 1500
$ws.Cells.Item(17, 2).Value = "falling"
$ws.Cells.Item(17, 3).Value = -1.988187789916992
$ws.Cells.Item(17, 4).Value = 1.761092185974121
$ws.Cells.Item(17, 5).Value = 0.5642168521881104
$ws.Cells.Item(17, 6).Value = -0.02540402729874076
$ws.Cells.Item(17, 7).Value = 0.01100221381563203
$ws.Cells.Item(17, 8).Value = 0.005378270172513961

$ws.Cells.Item(18, 1).Value = 1600
$ws.Cells.Item(18, 2).Value = "falling"
$ws.Cells.Item(18, 3).Value = -1.921907901763916
$ws.Cells.Item(18, 4).Value = 1.732572555541992
$ws.Cells.Item(18, 5).Value = 0.5020300149917603
$ws.Cells.Item(18, 6).Value = 0.0394008085131645
$ws.Cells.Item(18, 7).Value = -0.0215329993516206
$ws.Cells.Item(18, 8).Value = -0.001527163083665

$ws.Cells.Item(19, 1).Value = 1700
$ws.Cells.Item(19, 2).Value = "falling"
$ws.Cells.Item(19, 3).Value = -1.923232555389404
$ws.Cells.Item(19, 4).Value = 1.782273769378662
$ws.Cells.Item(19, 5).Value = 0.3185268342494964
$ws.Cells.Item(19, 6).Value = -0.005351710416700385
$ws.Cells.Item(19, 7).Value = -0.006526962327568362
$ws.Cells.Item(19, 8).Value = -0.02914889580468922

$ws.Cells.Item(20, 1).Value = 1800
$ws.Cells.Item(20, 2).Value = "falling"
$ws.Cells.Item(20, 3).Value = -1.910520076751709
$ws.Cells.Item(20, 4).Value = 1.814604759216309
$ws.Cells.Item(20, 5).Value = 0.3955313861370086
$ws.Cells.Item(20, 6).Value = 0.0001261571584187555
$ws.Cells.Item(20, 7).Value = 0.01954104799939229
$ws.Cells.Item(20, 8).Value = 0.05001127266365517

$ws.Cells.Item(21, 1).Value = 1900
$ws.Cells.Item(21, 2).Value = "falling"
$ws.Cells.Item(21, 3).Value = -2.009746074676514
$ws.Cells.Item(21, 4).Value = 1.783975601196289
$ws.Cells.Item(21, 5).Value = 0.4257155656814575
$ws.Cells.Item(21, 6).Value = 0.1110513180737254
$ws.Cells.Item(21, 7).Value = 0.04867666449559752
$ws.Cells.Item(21, 8).Value = -0.1655312067140687

$ws.Cells.Item(22, 1).Value = 2000
$ws.Cells.Item(22, 2).Value = "falling"
$ws.Cells.Item(22, 3).Value = -1.97176456451416
$ws.Cells.Item(22, 4).Value = 1.745009422302246
$ws.Cells.Item(22, 5).Value = 0.4838592410087585
$ws.Cells.Item(22, 6).Value = 0.1914996167887811
$ws.Cells.Item(22, 7).Value = 0.03030422819859344
$ws.Cells.Item(22, 8).Value = 0.02057685541069637

$ws.Cells.Item(23, 1).Value = 2100
$ws.Cells.Item(23, 2).Value = "falling"
$ws.Cells.Item(23, 3).Value = -1.845728397369385
$ws.Cells.Item(23, 4).Value = 1.672563552856445
$ws.Cells.Item(23, 5).Value = 0.5211508870124817
$ws.Cells.Item(23, 6).Value = 0.1882859338884768
$ws.Cells.Item(23, 7).Value = 0.0655750582480559
$ws.Cells.Item(23, 8).Value = 0.4253946024438608

$ws.Cells.Item(24, 1).Value = 2200
$ws.Cells.Item(24, 2).Value = "falling"
$ws.Cells.Item(24, 3).Value = -2.323234558105469
$ws.Cells.Item(24, 4).Value = 1.572612762451172
$ws.Cells.Item(24, 5).Value = 0.7384862899780273
$ws.Cells.Item(24, 6).Value = 0.06060181622919789
$ws.Cells.Item(24, 7).Value = 0.2589603908683945
$ws.Cells.Item(24, 8).Value = 0.4051430983387903

$ws.Cells.Item(25, 1).Value = 2300
$ws.Cells.Item(25, 2).Value = "falling"
$ws.Cells.Item(25, 3).Value = -3.269093036651612
$ws.Cells.Item(25, 4).Value = 0.9374399185180664
$ws.Cells.Item(25, 5).Value = 1.780844926834106
$ws.Cells.Item(25, 6).Value = -0.2129994086597282
$ws.Cells.Item(25, 7).Value = 0.3793273995751927
$ws.Cells.Item(25, 8).Value = 0.01590242139671223

$ws.Cells.Item(26, 1).Value = 2400
$ws.Cells.Item(26, 2).Value = "falling"
$ws.Cells.Item(26, 3).Value = -1.246297359466553
$ws.Cells.Item(26, 4).Value = 0.5547242164611816
$ws.Cells.Item(26, 5).Value = 1.987784385681152
$ws.Cells.Item(26, 6).Value = 0.5540348317312142
$ws.Cells.Item(26, 7).Value = -0.6639242237028846
$ws.Cells.Item(26, 8).Value = -1.121634874006974

$ws.Cells.Item(27, 1).Value = 2500
$ws.Cells.Item(27, 2).Value = "falling"
$ws.Cells.Item(27, 3).Value = -2.729167938232422
$ws.Cells.Item(27, 4).Value = 0.9446659088134766
$ws.Cells.Item(27, 5).Value = 2.494723320007324
$ws.Cells.Item(27, 6).Value = -0.05442012457743672
$ws.Cells.Item(27, 7).Value = -0.6851783446643664
$ws.Cells.Item(27, 8).Value = -1.266947746276855

$ws.Cells.Item(28, 1).Value = 2600
$ws.Cells.Item(28, 2).Value = "falling"
$ws.Cells.Item(28, 3).Value = -2.581250190734864
$ws.Cells.Item(28, 4).Value = -1.341118812561035
$ws.Cells.Item(28, 5).Value = 3.135838031768799
$ws.Cells.Item(28, 6).Value = -1.690450101442949
$ws.Cells.Item(28, 7).Value = 0.1374446801517266
$ws.Cells.Item(28, 8).Value = -2.446953586910074

$ws.Cells.Item(29, 1).Value = 2700
$ws.Cells.Item(29, 2).Value = "falling"
$ws.Cells.Item(29, 3).Value = -3.531517744064331
$ws.Cells.Item(29, 4).Value = -1.318747520446777
$ws.Cells.Item(29, 5).Value = 3.353370428085327
$ws.Cells.Item(29, 6).Value = -0.7898089613603607
$ws.Cells.Item(29, 7).Value = 0.3619974125986522
$ws.Cells.Item(29, 8).Value = -1.304489494017943

$ws.Cells.Item(30, 1).Value = 2800
$ws.Cells.Item(30, 2).Value = "falling"
$ws.Cells.Item(30, 3).Value = -4.357868194580078
$ws.Cells.Item(30, 4).Value = -0.8888802528381348
$ws.Cells.Item(30, 5).Value = 1.668703079223633
$ws.Cells.Item(30, 6).Value = 0.258422552243524
$ws.Cells.Item(30, 7).Value = 0.7682891954546356
$ws.Cells.Item(30, 8).Value = -0.810219791272407

$ws.Cells.Item(31, 1).Value = 2900
$ws.Cells.Item(31, 2).Value = "falling"
$ws.Cells.Item(31, 3).Value = -2.958048820495605
$ws.Cells.Item(31, 4).Value = 1.071966171264648
$ws.Cells.Item(31, 5).Value = 2.359493732452393
$ws.Cells.Item(31, 6).Value = 0.04960623909921762
$ws.Cells.Item(31, 7).Value = 0.7135171605193199
$ws.Cells.Item(31, 8).Value = -0.9448558195777682
